# Update the Max (column C) threshold values for rows 2-5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12.2
$ws.Range("C3").Value = 11.1
$ws.Range("C4").Value = 1.65
$ws.Range("C5").Value = 28

# The last row (row 6, a blank spacer row) is no longer needed - remove it
# so the used range shrinks back down to A1:C5.
$ws.Rows.Item(6).Delete()

# Widen columns A and C (Parameter name / Max) so the longer labels and
# decimal values are fully visible; column B goes back to its default width.
$ws.Columns.Item(1).ColumnWidth = 26.285714285714285
$ws.Columns.Item(3).ColumnWidth = 26.571428571428573

# Leave the selection on C3, matching where editing finished
$ws.Range("C3").Select()

# Set the page to print on Letter-ish A4 portrait (paper size 9 = A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
